# Update rows 2-9 on the active sheet to reflect the "Elemento Non Trovato (JS)"
# status recovered by the background worker, clearing the Provvedimento,
# Data Provvedimento and Protocollo uscita columns, and setting the note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$message = "Elemento Non Trovato (JS)"
$note = "Stato recuperato: Elemento Non Trovato (JS)."

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = $message   # C: Stato
    $ws.Cells.Item($row, 4).Value = $null      # D: Provvedimento
    $ws.Cells.Item($row, 5).Value = $null      # E: Data Provvedimento
    $ws.Cells.Item($row, 7).Value = $null      # G: Protocollo uscita
    $ws.Cells.Item($row, 8).Value = $note      # H: Note Usmaf
}
